$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ASR results data (rows 2-15) with new transcription values
$data = @(
    @{ Row = 2;  A = "<zero>";      B = "<kembe>";    C = 20 },
    @{ Row = 3;  A = "<part>";      B = "<part>";     C = 27 },
    @{ Row = 4;  A = "<zulu>";      B = "<zulu>";     C = 26 },
    @{ Row = 5;  A = "<water>";     B = "<long>";     C = 34 },
    @{ Row = 6;  A = "<can>";       B = "<can>";      C = 29 },
    @{ Row = 7;  A = "<a>";         B = "<up>";       C = 32 },
    @{ Row = 8;  A = "<you>";       B = "<you>";      C = 35 },
    @{ Row = 9;  A = "<number>";    B = "<none>";     C = 23 },
    @{ Row = 10; A = "<four>";      B = "<for>";      C = 28 },
    @{ Row = 11; A = "<word>";      B = "<would>";    C = 29 },
    @{ Row = 12; A = "<tango>";     B = "<tango>";    C = 32 },
    @{ Row = 13; A = "<backspace>"; B = "<backspace>"; C = 24 },
    @{ Row = 14; A = "<said>";      B = "<said>";     C = 29 },
    @{ Row = 15; A = "<november>";  B = "<cad>";      C = 7 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}

# Remove the now-unused trailing rows (16-18)
$ws.Rows("16:18").Delete()
